# Natmi following Dr Hou advice
# Rewrites the Il16-Cd4 ligand-receptor pair data: the underlying NATMI
# computation now uses 3 cluster replicates (E column = 3, previously 1) and
# expands the sender/target cluster combinations from 3 senders x 4 targets
# (rows 2-13) to the full 4 senders x 4 targets cartesian product (rows 2-17),
# with recomputed summary statistics in columns G-T.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array holds columns A..T (Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster, then the 16 numeric NATMI statistics).
$data = @(
    @("ECs", "Il16", "Cd4", "ECs", 3, 1, 6.566789, 19.700367, 0.3634074580963036, 0.3634074580963036, 1, 0.3333333333333333, 0.7633283333333334, 2.289985, 0.07247598971834317, 0.07247598971834318, 5.012616102721667, 45.113544924495, 0.02633831519655692, 0.02633831519655693),
    @("ECs", "Il16", "Cd4", "FAPs", 3, 1, 6.566789, 19.700367, 0.3634074580963036, 0.3634074580963036, 3, 1, 4.451606333333333, 13.354819, 0.4226681504613934, 0.4226681504613934, 29.23275950206367, 263.094835518573, 0.153600758177441, 0.153600758177441),
    @("ECs", "Il16", "Cd4", "M2", 3, 1, 6.566789, 19.700367, 0.3634074580963036, 0.3634074580963036, 3, 1, 4.617562, 13.852686, 0.438425198465246, 0.438425198465246, 30.322555348418, 272.902998135762, 0.1593269869396225, 0.1593269869396225),
    @("ECs", "Il16", "Cd4", "sCs", 3, 1, 6.566789, 19.700367, 0.3634074580963036, 0.3634074580963036, 3, 1, 0.699658, 2.098974, 0.06643066135501745, 0.06643066135501746, 4.594506458162, 41.350558123458, 0.02414139778268324, 0.02414139778268324),
    @("FAPs", "Il16", "Cd4", "ECs", 3, 1, 6.440526333333334, 19.321579, 0.3564200560729107, 0.3564200560729107, 1, 0.3333333333333333, 0.7633283333333334, 2.289985, 0.07247598971834317, 0.07247598971834318, 4.916236231812778, 44.246126086315, 0.02583189631935157, 0.02583189631935157),
    @("FAPs", "Il16", "Cd4", "FAPs", 3, 1, 6.440526333333334, 19.321579, 0.3564200560729107, 0.3564200560729107, 3, 1, 4.451606333333333, 13.354819, 0.4226681504613934, 0.4226681504613934, 28.67068781546678, 258.036190339201, 0.1506474058876833, 0.1506474058876833),
    @("FAPs", "Il16", "Cd4", "M2", 3, 1, 6.440526333333334, 19.321579, 0.3564200560729107, 0.3564200560729107, 3, 1, 4.617562, 13.852686, 0.438425198465246, 0.438425198465246, 29.73952965679934, 267.655766911194, 0.15626353382076, 0.15626353382076),
    @("FAPs", "Il16", "Cd4", "sCs", 3, 1, 6.440526333333334, 19.321579, 0.3564200560729107, 0.3564200560729107, 3, 1, 0.699658, 2.098974, 0.06643066135501745, 0.06643066135501746, 4.506165773327334, 40.555491959946, 0.02367722004511586, 0.02367722004511586),
    @("M2", "Il16", "Cd4", "ECs", 3, 1, 4.846355333333333, 14.539066, 0.2681983040292799, 0.2681983040292799, 1, 0.3333333333333333, 0.7633283333333334, 2.289985, 0.07247598971834317, 0.07247598971834318, 3.699360339334445, 33.29424305401, 0.01943793752530316, 0.01943793752530316),
    @("M2", "Il16", "Cd4", "FAPs", 3, 1, 4.846355333333333, 14.539066, 0.2681983040292799, 0.2681983040292799, 3, 1, 4.451606333333333, 13.354819, 0.4226681504613934, 0.4226681504613934, 21.57406609545044, 194.166594859054, 0.1133588811209382, 0.1133588811209382),
    @("M2", "Il16", "Cd4", "M2", 3, 1, 4.846355333333333, 14.539066, 0.2681983040292799, 0.2681983040292799, 3, 1, 4.617562, 13.852686, 0.438425198465246, 0.438425198465246, 22.37834622569734, 201.405116031276, 0.1175848946720794, 0.1175848946720794),
    @("M2", "Il16", "Cd4", "sCs", 3, 1, 4.846355333333333, 14.539066, 0.2681983040292799, 0.2681983040292799, 3, 1, 0.699658, 2.098974, 0.06643066135501745, 0.06643066135501746, 3.390791279809333, 30.517121518284, 0.0178165907109591, 0.0178165907109591),
    @("sCs", "Il16", "Cd4", "ECs", 3, 1, 0.216374, 0.6491220000000001, 0.01197418180150597, 0.01197418180150597, 1, 0.3333333333333333, 0.7633283333333334, 2.289985, 0.07247598971834317, 0.07247598971834318, 0.1651644047966667, 1.48647964317, 0.0008678406771315188, 0.0008678406771315189),
    @("sCs", "Il16", "Cd4", "FAPs", 3, 1, 0.216374, 0.6491220000000001, 0.01197418180150597, 0.01197418180150597, 3, 1, 4.451606333333333, 13.354819, 0.4226681504613934, 0.4226681504613934, 0.9632118687686668, 8.668906818918, 0.005061105275331005, 0.005061105275331004),
    @("sCs", "Il16", "Cd4", "M2", 3, 1, 0.216374, 0.6491220000000001, 0.01197418180150597, 0.01197418180150597, 3, 1, 4.617562, 13.852686, 0.438425198465246, 0.438425198465246, 0.9991203601880002, 8.992083241692001, 0.005249783032784194, 0.005249783032784193),
    @("sCs", "Il16", "Cd4", "sCs", 3, 1, 0.216374, 0.6491220000000001, 0.01197418180150597, 0.01197418180150597, 3, 1, 0.699658, 2.098974, 0.06643066135501745, 0.06643066135501746, 0.151387800092, 1.362490200828, 0.000795452816259256, 0.0007954528162592561)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowValues[$c]
    }
}
